# Generate Report for Handback
#
# The handback generator re-ran and produced a fresh localization-status
# report:
#   - Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is shown (the
#     Overview rollup columns for zh-cn/de-de, and each locale sheet's own
#     Status column).
#   - Each locale sheet's "Latest Handback DateTime" is refreshed to the
#     moment this handback report was generated.
#   - The stale "handback file is not the latest" warning that used to sit
#     in "Error Detail" is cleared now that the handback is in sync, for
#     both locale sheets.
#   - The Status / Error Detail columns are re-sized to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status roll-up columns for zh-cn (E) and de-de (F) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn sheet: Status, Latest Handback DateTime, Error Detail ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-17 14:49:56"
$zhcn.Range("P2").Value = ""

# --- de-de sheet: Status, Latest Handback DateTime, Error Detail ---
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-17 14:50:22"
$dede.Range("P2").Value = ""

# --- Resize the Status / Error Detail columns to match the new text ---
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

$zhcn.Columns.Item(3).ColumnWidth  = 29.144371396019366
$zhcn.Columns.Item(16).ColumnWidth = 12.913719540550566

$dede.Columns.Item(3).ColumnWidth  = 29.144371396019366
$dede.Columns.Item(16).ColumnWidth = 12.913719540550566
